$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25, shifting the existing
# rows 25-40 down to 26-41 (this also pushes their formatting down,
# which keeps the date style on column D intact).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44813
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100114007
$ws.Range("G25").Value = "Jengibre"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13400
$ws.Range("N25").Value = "$/caja 13 kilos"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 1031
$ws.Range("Q25").Value = 13
$ws.Range("R25").Value = "Hortaliza"
